# Update NATMI LR-pair (Cd200-Cd200r1) sheet with re-computed TPM-based values.
# Ligand-side stats (G,H,I,J) are keyed by the "Sending cluster" (column A).
# Receptor-side stats (K,L,M,N,O,P) are keyed by the "Target cluster" (column D).
# Edge weights (Q,R,S,T) are simply the products of the matching ligand/receptor
# average & total expression (and their specificity) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New per-cluster ligand expression stats: avg, total, avg-specificity, total-specificity
$ligandStats = @{
    "ECs"              = @(38.84293433333333, 116.528803, 0.4759825783603507, 0.5627959792472441)
    "FAPs"             = @(2.977970666666666, 8.933911999999999, 0.03649214922944397, 0.043147870939245)
    "Inflammatory-Mac" = @(1.675036, 5.025107999999999, 0.02052594552420851, 0.02426962694951188)
    "MuSCs"            = @(37.764002, 75.52800400000001, 0.4627613065200399, 0.3647755394155193)
    "Resolving-Mac"    = @(0.345847, 1.037541, 0.004238020365956876, 0.005010983448479815)
}

# New per-cluster receptor expression stats: cells, detection-rate, avg, total, avg-specificity, total-specificity
$receptorStats = @{
    "Inflammatory-Mac" = @(3, 1, 8.366847333333332, 25.100542, 0.3599008094639652, 0.3608281410662089)
    "MuSCs"            = @(2, 1, 0.17924, 0.35848, 0.007710027268134825, 0.005153262109217186)
    "Resolving-Mac"    = @(3, 1, 14.70156066666667, 44.104682, 0.6323891632679, 0.6340185968245738)
}

$lastRow = $ws.Range("A1").End(4).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $sendCluster = $ws.Cells.Item($r, 1).Value2
    $targetCluster = $ws.Cells.Item($r, 4).Value2

    $lig = $ligandStats[$sendCluster]
    $rec = $receptorStats[$targetCluster]

    $ligAvg = $lig[0]
    $ligTot = $lig[1]
    $ligAvgSpec = $lig[2]
    $ligTotSpec = $lig[3]

    $recCells = $rec[0]
    $recRate = $rec[1]
    $recAvg = $rec[2]
    $recTot = $rec[3]
    $recAvgSpec = $rec[4]
    $recTotSpec = $rec[5]

    $ws.Cells.Item($r, 7).Value = $ligAvg
    $ws.Cells.Item($r, 8).Value = $ligTot
    $ws.Cells.Item($r, 9).Value = $ligAvgSpec
    $ws.Cells.Item($r, 10).Value = $ligTotSpec

    $ws.Cells.Item($r, 11).Value = $recCells
    $ws.Cells.Item($r, 12).Value = $recRate
    $ws.Cells.Item($r, 13).Value = $recAvg
    $ws.Cells.Item($r, 14).Value = $recTot
    $ws.Cells.Item($r, 15).Value = $recAvgSpec
    $ws.Cells.Item($r, 16).Value = $recTotSpec

    $ws.Cells.Item($r, 17).Value = $ligAvg * $recAvg
    $ws.Cells.Item($r, 18).Value = $ligTot * $recTot
    $ws.Cells.Item($r, 19).Value = $ligAvgSpec * $recAvgSpec
    $ws.Cells.Item($r, 20).Value = $ligTotSpec * $recTotSpec
}
